# Diverse småoppdateringer i ekstra
# Adds a new row for "2021 - Vår" below the existing exam-archive table,
# and updates the selected cell to reflect the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "2021 - Vår"
$ws.Range("B10").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/hjemme-21-v.pdf)"
$ws.Range("C10").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/hjemme-21-v-solprop.pdf)"
$ws.Range("D10").Value = "Data ikke lenger tilgjengelig."

$ws.Range("C11").Select()
